$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2.. down by one)
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the base_dir setting
$ws.Cells.Item(2, 1).Value = "base_dir"
$ws.Cells.Item(2, 2).Value = "os.getcwd()"

# Put selection on B5 to mirror the saved workbook view state
$ws.Range("B5").Select()
